$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold plain text that looks numeric (e.g. "1.003", "99.75").
# Assigning .Value directly would make Excel auto-convert these to real numbers
# (losing the original text formatting / trailing zeros), so we temporarily mark
# each target cell as Text, write the literal string, then clear the temporary
# number format again so the cell style matches the original (unformatted) cell.
$textRows = @(2, 3, 4, 5, 7, 8, 9, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 29, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51)
foreach ($r in $textRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "24.684.60"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "1.692.88"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "316.81"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "0.3949"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "0.4057"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").Value = "1.488"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("D11").Value = "52.26"
$ws.Range("E11").Value = "  -2.82%  "
$ws.Range("D12").Value = "0.08875"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").Value = "7.255"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "23.62"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").Value = "8.056"
$ws.Range("E15").Value = "  +7.39%  "
$ws.Range("E16").Value = "  +0.14%  "
$ws.Range("D17").Value = "1.696.78"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "99.75"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("D19").Value = "0.07026"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("D20").Value = "19.60"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "7.007"
$ws.Range("E21").Value = "  +4.92%  "
$ws.Range("D22").Value = "1.006"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "14.36"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").Value = "24.674.40"
$ws.Range("E24").Value = "  +1.09%  "
$ws.Range("D25").Value = "3.206"
$ws.Range("E25").Value = "  +7.21%  "
$ws.Range("D26").Value = "2.360"
$ws.Range("E26").Value = "  +1.84%  "
$ws.Range("D27").Value = "22.74"
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").Value = "162.30"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("D29").Value = "135.88"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("D31").Value = "7.596"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").Value = "1.883.05"
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "0.08613"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "1.057"
$ws.Range("E34").Value = "  -2.14%  "
$ws.Range("D35").Value = "7.101"
$ws.Range("E35").Value = "  -3.15%  "
$ws.Range("D36").Value = "11.33"
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").Value = "0.2730"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").Value = "1.890"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("D39").Value = "14.49"
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").Value = "0.09212"
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("D41").Value = "0.02725"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "1.471"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").Value = "0.7665"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "16.00"
$ws.Range("E44").Value = "  +3.45%  "
$ws.Range("D45").Value = "2.597"
$ws.Range("E45").Value = "  +6.62%  "
$ws.Range("D46").Value = "0.7156"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").Value = "4.226"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "140.08"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "1.319"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").Value = "90.61"
$ws.Range("E51").Value = "  +5.04%  "

foreach ($r in $textRows) {
    $ws.Range("D$r").ClearFormats()
}

Write-Output "done"
